# "Group reflection" slide (slide 8) — update the group's answers.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Paragraph 2 (level-1 bullet under "How did you feel about this milestone? ..."):
# "It was good, smartdraw sucks" -> "We became more comfortable using smartdraw
# which made making the future state diagram easyer"
$para = $tr.Paragraphs(2, 1)
$para.Text = "We became more comfortable using smartdraw which made making the future state diagram easyer"

# Paragraph 7 (level-1 bullet under "How will you use what you have learned going forward?"):
# "We will use our future state diagram for milestone 3" ->
# "We will use our future state diagram for our user story map in milestone 3"
$para = $tr.Paragraphs(7, 1)
$para.Text = "We will use our future state diagram for our user story map in milestone 3"

# Paragraph 9 (level-1 bullet under the "stuff & things" question):
# "Stick people!!!!!!!" -> "We are not sure what to put for participant name"
$para = $tr.Paragraphs(9, 1)
$para.Text = "We are not sure what to put for participant name"
